$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

# Row 242: event_id 14316326
Set-TextCell 242 1 '14316326'
Set-TextCell 242 2 '2025-08-18'
Set-TextCell 242 3 'Jannik Sinner'
Set-TextCell 242 4 'Carlos Alcaraz'
Set-TextCell 242 5 'Gana Carlos Alcaraz'
$ws.Cells.Item(242, 6).Value = 2.5
Set-TextCell 242 7 ""
Set-TextCell 242 8 ""

# Row 243: event_id 14452482
Set-TextCell 243 1 '14452482'
Set-TextCell 243 2 '2025-08-18'
Set-TextCell 243 3 'Hugo Gaston'
Set-TextCell 243 4 'Marton Fucsovics'
Set-TextCell 243 5 'Gana Hugo Gaston'
$ws.Cells.Item(243, 6).Value = 2.75
Set-TextCell 243 7 ""
Set-TextCell 243 8 ""

# Row 244: event_id 14466941
Set-TextCell 244 1 '14466941'
Set-TextCell 244 2 '2025-08-18'
Set-TextCell 244 3 'Raphael Collignon'
Set-TextCell 244 4 'Nishesh Basavareddy'
Set-TextCell 244 5 'Gana Raphael Collignon'
$ws.Cells.Item(244, 6).Value = 3.75
Set-TextCell 244 7 ""
Set-TextCell 244 8 ""

# Row 245: event_id 14452485
Set-TextCell 245 1 '14452485'
Set-TextCell 245 2 '2025-08-18'
Set-TextCell 245 3 'Roman Safiullin'
Set-TextCell 245 4 'Christopher O''Connell'
Set-TextCell 245 5 'Gana Christopher O''Connell'
$ws.Cells.Item(245, 6).Value = 2.63
Set-TextCell 245 7 ""
Set-TextCell 245 8 ""

# Row 246: event_id 14428726
Set-TextCell 246 1 '14428726'
Set-TextCell 246 2 '2025-08-19'
Set-TextCell 246 3 'Nuno Borges'
Set-TextCell 246 4 'Kamil Majchrzak'
Set-TextCell 246 5 'Gana Nuno Borges'
$ws.Cells.Item(246, 6).Value = 1.73
Set-TextCell 246 7 ""
Set-TextCell 246 8 ""

# Row 247: event_id 14316411
Set-TextCell 247 1 '14316411'
Set-TextCell 247 2 '2025-08-18'
Set-TextCell 247 3 'Iga Swiatek'
Set-TextCell 247 4 'Jasmine Paolini'
Set-TextCell 247 5 'Gana Jasmine Paolini'
$ws.Cells.Item(247, 6).Value = 5
Set-TextCell 247 7 ""
Set-TextCell 247 8 ""

# Row 248: event_id 14460690
Set-TextCell 248 1 '14460690'
Set-TextCell 248 2 '2025-08-18'
Set-TextCell 248 3 'Marie Bouzkova'
Set-TextCell 248 4 'Zeynep Sonmez'
Set-TextCell 248 5 'Gana Zeynep Sonmez'
$ws.Cells.Item(248, 6).Value = 4.33
Set-TextCell 248 7 ""
Set-TextCell 248 8 ""

# Row 249: event_id 14460694
Set-TextCell 249 1 '14460694'
Set-TextCell 249 2 '2025-08-19'
Set-TextCell 249 3 'Jaqueline Cristian'
Set-TextCell 249 4 'Leylah Fernandez'
Set-TextCell 249 5 'Gana Jaqueline Cristian'
$ws.Cells.Item(249, 6).Value = 3
Set-TextCell 249 7 ""
Set-TextCell 249 8 ""

# Row 250: event_id 14452725
Set-TextCell 250 1 '14452725'
Set-TextCell 250 2 '2025-08-18'
Set-TextCell 250 3 'Elena-Gabriela Ruse'
Set-TextCell 250 4 'Hailey Baptiste'
Set-TextCell 250 5 'Gana Hailey Baptiste'
$ws.Cells.Item(250, 6).Value = 1.91
Set-TextCell 250 7 ""
Set-TextCell 250 8 ""

# Row 251: event_id 14466796
Set-TextCell 251 1 '14466796'
Set-TextCell 251 2 '2025-08-18'
Set-TextCell 251 3 'Sorana Cirstea'
Set-TextCell 251 4 'Moyuka Uchijima'
Set-TextCell 251 5 'Gana Moyuka Uchijima'
$ws.Cells.Item(251, 6).Value = 2.75
Set-TextCell 251 7 ""
Set-TextCell 251 8 ""

# Row 252: event_id 14466797
Set-TextCell 252 1 '14466797'
Set-TextCell 252 2 '2025-08-18'
Set-TextCell 252 3 'Talia Gibson'
Set-TextCell 252 4 'Greet Minnen'
Set-TextCell 252 5 'Gana Talia Gibson'
$ws.Cells.Item(252, 6).Value = 2.63
Set-TextCell 252 7 ""
Set-TextCell 252 8 ""

# Row 253: event_id 14452707
Set-TextCell 253 1 '14452707'
Set-TextCell 253 2 '2025-08-18'
Set-TextCell 253 3 'Katie Boulter'
Set-TextCell 253 4 'Yue Yuan'
Set-TextCell 253 5 'Gana Yue Yuan'
$ws.Cells.Item(253, 6).Value = 2
Set-TextCell 253 7 ""
Set-TextCell 253 8 ""

# Row 254: event_id 14466731
Set-TextCell 254 1 '14466731'
Set-TextCell 254 2 '2025-08-18'
Set-TextCell 254 3 'Andres Martin'
Set-TextCell 254 4 'Viktor Durasovic'
Set-TextCell 254 5 'Gana Viktor Durasovic'
$ws.Cells.Item(254, 6).Value = 3.5
Set-TextCell 254 7 ""
Set-TextCell 254 8 ""

# Row 255: event_id 14466726
Set-TextCell 255 1 '14466726'
Set-TextCell 255 2 '2025-08-18'
Set-TextCell 255 3 'Clement Tabur'
Set-TextCell 255 4 'Martin Landaluce'
Set-TextCell 255 5 'Gana Clement Tabur'
$ws.Cells.Item(255, 6).Value = 3.4
Set-TextCell 255 7 ""
Set-TextCell 255 8 ""

# Row 256: event_id 14466727
Set-TextCell 256 1 '14466727'
Set-TextCell 256 2 '2025-08-18'
Set-TextCell 256 3 'Hugo Grenier'
Set-TextCell 256 4 'Omar Jasika'
Set-TextCell 256 5 'Gana Omar Jasika'
$ws.Cells.Item(256, 6).Value = 3.4
Set-TextCell 256 7 ""
Set-TextCell 256 8 ""

# Row 257: event_id 14466730
Set-TextCell 257 1 '14466730'
Set-TextCell 257 2 '2025-08-18'
Set-TextCell 257 3 'Pablo Llamas Ruiz'
Set-TextCell 257 4 'Yannick Hanfmann'
Set-TextCell 257 5 'Gana Pablo Llamas Ruiz'
$ws.Cells.Item(257, 6).Value = 2.75
Set-TextCell 257 7 ""
Set-TextCell 257 8 ""

# Row 258: event_id 14466733
Set-TextCell 258 1 '14466733'
Set-TextCell 258 2 '2025-08-18'
Set-TextCell 258 3 'Cristian Garin'
Set-TextCell 258 4 'Yosuke Watanuki'
Set-TextCell 258 5 'Gana Yosuke Watanuki'
$ws.Cells.Item(258, 6).Value = 2.2
Set-TextCell 258 7 ""
Set-TextCell 258 8 ""

# Row 259: event_id 14466778
Set-TextCell 259 1 '14466778'
Set-TextCell 259 2 '2025-08-18'
Set-TextCell 259 3 'Jesper De Jong'
Set-TextCell 259 4 'Vilius Gaubas'
Set-TextCell 259 5 'Gana Vilius Gaubas'
$ws.Cells.Item(259, 6).Value = 4.33
Set-TextCell 259 7 ""
Set-TextCell 259 8 ""

# Row 260: event_id 14466729
Set-TextCell 260 1 '14466729'
Set-TextCell 260 2 '2025-08-18'
Set-TextCell 260 3 'Lukas Klein'
Set-TextCell 260 4 'Federico Agustin Gomez'
Set-TextCell 260 5 'Gana Federico Agustin Gomez'
$ws.Cells.Item(260, 6).Value = 3.4
Set-TextCell 260 7 ""
Set-TextCell 260 8 ""

# Row 261: event_id 14466718
Set-TextCell 261 1 '14466718'
Set-TextCell 261 2 '2025-08-18'
Set-TextCell 261 3 'Mikhail Kukushkin'
Set-TextCell 261 4 'Andrea Pellegrino'
Set-TextCell 261 5 'Gana Andrea Pellegrino'
$ws.Cells.Item(261, 6).Value = 2.5
Set-TextCell 261 7 ""
Set-TextCell 261 8 ""

# Row 262: event_id 14466779
Set-TextCell 262 1 '14466779'
Set-TextCell 262 2 '2025-08-18'
Set-TextCell 262 3 'Clement Chidekh'
Set-TextCell 262 4 'Jan-Lennard Struff'
Set-TextCell 262 5 'Gana Clement Chidekh'
$ws.Cells.Item(262, 6).Value = 4.5
Set-TextCell 262 7 ""
Set-TextCell 262 8 ""

# Row 263: event_id 14466754
Set-TextCell 263 1 '14466754'
Set-TextCell 263 2 '2025-08-18'
Set-TextCell 263 3 'Alexis Galarneau'
Set-TextCell 263 4 'Dino Prižmić'
Set-TextCell 263 5 'Gana Alexis Galarneau'
$ws.Cells.Item(263, 6).Value = 4.33
Set-TextCell 263 7 ""
Set-TextCell 263 8 ""

# Row 264: event_id 14466745
Set-TextCell 264 1 '14466745'
Set-TextCell 264 2 '2025-08-18'
Set-TextCell 264 3 'Francesco Passaro'
Set-TextCell 264 4 'Yu Hsiou Hsu'
Set-TextCell 264 5 'Gana Francesco Passaro'
$ws.Cells.Item(264, 6).Value = 1.83
Set-TextCell 264 7 ""
Set-TextCell 264 8 ""

# Row 265: event_id 14466751
Set-TextCell 265 1 '14466751'
Set-TextCell 265 2 '2025-08-18'
Set-TextCell 265 3 'Pol Martin Tiffon'
Set-TextCell 265 4 'Fajing Sun'
Set-TextCell 265 5 'Gana Pol Martin Tiffon'
$ws.Cells.Item(265, 6).Value = 2.25
Set-TextCell 265 7 ""
Set-TextCell 265 8 ""

# Row 266: event_id 14466843
Set-TextCell 266 1 '14466843'
Set-TextCell 266 2 '2025-08-18'
Set-TextCell 266 3 'Dalma Galfi'
Set-TextCell 266 4 'Akasha Urhobo'
Set-TextCell 266 5 'Gana Akasha Urhobo'
$ws.Cells.Item(266, 6).Value = 4.33
Set-TextCell 266 7 ""
Set-TextCell 266 8 ""

# Row 267: event_id 14466841
Set-TextCell 267 1 '14466841'
Set-TextCell 267 2 '2025-08-18'
Set-TextCell 267 3 'Sada Nahimana'
Set-TextCell 267 4 'Manon Leonard'
Set-TextCell 267 5 'Gana Sada Nahimana'
$ws.Cells.Item(267, 6).Value = 2.38
Set-TextCell 267 7 ""
Set-TextCell 267 8 ""

# Row 268: event_id 14466851
Set-TextCell 268 1 '14466851'
Set-TextCell 268 2 '2025-08-18'
Set-TextCell 268 3 'Tereza Valentova'
Set-TextCell 268 4 'Mona Barthel'
Set-TextCell 268 5 'Gana Mona Barthel'
$ws.Cells.Item(268, 6).Value = 4.33
Set-TextCell 268 7 ""
Set-TextCell 268 8 ""

# Row 269: event_id 14466817
Set-TextCell 269 1 '14466817'
Set-TextCell 269 2 '2025-08-18'
Set-TextCell 269 3 'Astra Sharma'
Set-TextCell 269 4 'Lola Radivojevic'
Set-TextCell 269 5 'Gana Lola Radivojevic'
$ws.Cells.Item(269, 6).Value = 2.25
Set-TextCell 269 7 ""
Set-TextCell 269 8 ""

# Row 270: event_id 14466828
Set-TextCell 270 1 '14466828'
Set-TextCell 270 2 '2025-08-18'
Set-TextCell 270 3 'Jana Fett'
Set-TextCell 270 4 'Joanna Garland'
Set-TextCell 270 5 'Gana Jana Fett'
$ws.Cells.Item(270, 6).Value = 2.25
Set-TextCell 270 7 ""
Set-TextCell 270 8 ""

# Row 271: event_id 14466831
Set-TextCell 271 1 '14466831'
Set-TextCell 271 2 '2025-08-18'
Set-TextCell 271 3 'Katie Volynets'
Set-TextCell 271 4 'Sofia Costoulas'
Set-TextCell 271 5 'Gana Sofia Costoulas'
$ws.Cells.Item(271, 6).Value = 3
Set-TextCell 271 7 ""
Set-TextCell 271 8 ""

# Row 272: event_id 14466833
Set-TextCell 272 1 '14466833'
Set-TextCell 272 2 '2025-08-18'
Set-TextCell 272 3 'Barbora Palicova'
Set-TextCell 272 4 'Elena Pridankina'
Set-TextCell 272 5 'Gana Barbora Palicova'
$ws.Cells.Item(272, 6).Value = 2.2
Set-TextCell 272 7 ""
Set-TextCell 272 8 ""

# Row 273: event_id 14466819
Set-TextCell 273 1 '14466819'
Set-TextCell 273 2 '2025-08-18'
Set-TextCell 273 3 'Bernarda Pera'
Set-TextCell 273 4 'Emerson Jones'
Set-TextCell 273 5 'Gana Emerson Jones'
$ws.Cells.Item(273, 6).Value = 2.38
Set-TextCell 273 7 ""
Set-TextCell 273 8 ""

# Row 274: event_id 14466815
Set-TextCell 274 1 '14466815'
Set-TextCell 274 2 '2025-08-18'
Set-TextCell 274 3 'Petra Marčinko'
Set-TextCell 274 4 'Xinyu Gao'
Set-TextCell 274 5 'Gana Xinyu Gao'
$ws.Cells.Item(274, 6).Value = 3
Set-TextCell 274 7 ""
Set-TextCell 274 8 ""

# Row 275: event_id 14466787
Set-TextCell 275 1 '14466787'
Set-TextCell 275 2 '2025-08-18'
Set-TextCell 275 3 'Daria Snigur'
Set-TextCell 275 4 'Simona Waltert'
Set-TextCell 275 5 'Gana Simona Waltert'
$ws.Cells.Item(275, 6).Value = 2.75
Set-TextCell 275 7 ""
Set-TextCell 275 8 ""

# Row 276: event_id 14466788
Set-TextCell 276 1 '14466788'
Set-TextCell 276 2 '2025-08-18'
Set-TextCell 276 3 'Lauren Davis'
Set-TextCell 276 4 'Hina Inoue'
Set-TextCell 276 5 'Gana Hina Inoue'
$ws.Cells.Item(276, 6).Value = 2.62
Set-TextCell 276 7 ""
Set-TextCell 276 8 ""

# Row 277: event_id 14458875
Set-TextCell 277 1 '14458875'
Set-TextCell 277 2 '2025-08-18'
Set-TextCell 277 3 'Alexander Ritschard'
Set-TextCell 277 4 'Rudolf Molleker'
Set-TextCell 277 5 'Gana Alexander Ritschard'
$ws.Cells.Item(277, 6).Value = 2
Set-TextCell 277 7 ""
Set-TextCell 277 8 ""
